$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits after
#    "Visit GitHub Website" (it will be re-added at the very end of the
#    document, next to the new "Thanks." paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) At the end of the document, after the final (empty, bold) paragraph,
#    add two more empty bold paragraphs followed by a bold "Thanks."
#    paragraph that carries the (re-created) "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xmlFragment = @"
<w:p $ns><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p $ns><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p $ns><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Thanks.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

$insertionPoint.InsertXML($xmlFragment)
